# fix(publipostage): Correct status name
#
# Updates the shared-string labels used in the "statut_label" (col B) and
# "statut_name" (col C) columns:
#   - "bleu" -> "noir"
#   - "résultat et / ou publication posté" -> "résultat postés ou publiés"
#   - "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
#   - "résultat et / ou publication posté dans les 12 mois"
#       -> "résultat postés ou publiés dans les 12 mois"
#   - "résultat et / ou publication posté dans les 36 mois"
#       -> "résultat postés ou publiés dans les 36 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

# Replace the two "... dans les N mois" variants first so the shorter base
# phrase below doesn't need to special-case them (Replace matches the
# substring wherever it occurs within a cell's text).
$cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés")
$cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
$cells.Replace("bleu", "noir")
